$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. "Délégation de la Recherche Clinique " -> split into 3 runs
# -----------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Délégation de la Recherche Clinique ")
if ($found) {
    $target = $d.Range($rng.Start, $rng.End)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r><w:t xml:space="preserve">Délégation </w:t></w:r>' +
           '<w:r><w:t>à</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> la Recherche Clinique </w:t></w:r>' +
           '</w:p>'
    $target.InsertXML($xml)
}

# -----------------------------------------------------------------
# 2. " du développement" -> " à l'Innovation"
# -----------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(" du développement")
if ($found) {
    $target = $d.Range($rng.Start, $rng.End)
    $target.Text = " à l’Innovation"
}

# -----------------------------------------------------------------
# 3. DATE_GEL run split: "DATE_GEL" -> "D" + "ATE_GE"(bookmarked) + "L"
# -----------------------------------------------------------------
$bmGel = $d.Bookmarks("DATE_GEL")
$s = $bmGel.Start
$e = $bmGel.End
$bmGel.Delete()
$sub = $d.Range($s + 1, $e - 1)
$d.Bookmarks.Add("DATE_GEL", $sub) | Out-Null

# -----------------------------------------------------------------
# 4. DATE_MAJ run split: "DATE_MAJ" -> "D" + "ATE_MA"(DATE_MAJ + _GoBack) + "J"
# -----------------------------------------------------------------
$bmMaj = $d.Bookmarks("DATE_MAJ")
$s2 = $bmMaj.Start
$e2 = $bmMaj.End
$bmMaj.Delete()
$sub2 = $d.Range($s2 + 1, $e2 - 1)
$d.Bookmarks.Add("DATE_MAJ", $sub2) | Out-Null
$d.Bookmarks.Add("_GoBack", $sub2) | Out-Null

# -----------------------------------------------------------------
# 5. Insert a section break (next page) right after the final
#    page-break paragraph (this also relocates _GoBack previously
#    at this spot - already handled above).
# -----------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Host "done"
